$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=0; B="Compass Error Compass Error Compass data error Please contact DJI Support ."; C="Compass Error"; D="0-1"; E="Event"; F="Event"},
    @{Row=3;  A=0; B="Compass Error Compass Error Compass data error Please contact DJI Support ."; C="Please contact DJI Support"; D="7-10"; E="NonEvent"; F="NonEvent"},
    @{Row=4;  A=0; B="Compass Error Compass Error Compass data error Please contact DJI Support ."; C="Compass data error"; D="4-6"; E="Event"; F="Event"},
    @{Row=5;  A=0; B="Compass Error Compass Error Compass data error Please contact DJI Support ."; C="Compass Error"; D="2-3"; E="Event"; F="Event"},
    @{Row=6;  A=1; B="GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn ."; C="GPS signal weak"; D="0-2"; E="Event"; F="Event"},
    @{Row=7;  A=1; B="GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn ."; C="Aircraft in Altitude Zone"; D="6-9"; E="NonEvent"; F="NonEvent"},
    @{Row=8;  A=1; B="GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn ."; C="Fly with caution"; D="3-5"; E="NonEvent"; F="NonEvent"},
    @{Row=9;  A=1; B="GPS signal weak Fly with caution Aircraft in Altitude Zone Max altitude set to nnn ."; C="Max altitude set to nnn"; D="10-14"; E="Event"; F="Event"},
    @{Row=10; A=2; B="High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP ."; C="High wind velocity"; D="0-2"; E="Event"; F="Event"},
    @{Row=11; A=2; B="High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP ."; C="Ensure the aircraft remains within your line of sight and fly with caution"; D="3-15"; E="NonEvent"; F="NonEvent"},
    @{Row=12; A=2; B="High wind velocity Ensure the aircraft remains within your line of sight and fly with caution High wind velocity Fly with caution and land in a safe place ASAP ."; C="Fly with caution and land in a safe place ASAP"; D="19-28"; E="NonEvent"; F="NonEvent"},
    @{Row=13; A=3; B="Motor speed error Land or return to home promptly After powering off the aircraft, replace the propeller on the beeping ESC If the issue persists, contact DJI Support ."; C="Land or return to home promptly"; D="3-8"; E="NonEvent"; F="NonEvent"},
    @{Row=14; A=3; B="Motor speed error Land or return to home promptly After powering off the aircraft, replace the propeller on the beeping ESC If the issue persists, contact DJI Support ."; C="Motor speed error"; D="0-2"; E="Event"; F="Event"},
    @{Row=15; A=3; B="Motor speed error Land or return to home promptly After powering off the aircraft, replace the propeller on the beeping ESC If the issue persists, contact DJI Support ."; C="If the issue persists, contact DJI Support"; D="21-27"; E="NonEvent"; F="NonEvent"},
    @{Row=16; A=3; B="Motor speed error Land or return to home promptly After powering off the aircraft, replace the propeller on the beeping ESC If the issue persists, contact DJI Support ."; C="After powering off the aircraft, replace the propeller on the beeping ESC"; D="9-20"; E="NonEvent"; F="NonEvent"}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
}
